# Remove erroneous / duplicate student rows that were picked up while
# filtering the roster by "Hoc luc" (Academic performance) & "Hanh kiem"
# (Conduct): student 2102030089 (Le Ho Nguyen Khoi) had been entered twice
# -- once on sheet "12t" and once on sheet "12i" -- and student 2102030073
# (Trao Le Hong Anh) on sheet "12i" was a stray duplicate row as well.
# Delete those rows (column B holds "Ma hoc sinh") wherever they are found.

$wb = $excel.ActiveWorkbook

function Remove-StudentRows($ws, $ids) {
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = $lastRow; $r -ge 2; $r--) {
        $mshs = $ws.Cells.Item($r, 2).Value()
        if ($mshs -ne $null -and ($ids -contains [string]$mshs)) {
            $ws.Rows.Item($r).Delete()
        }
    }
}

$ws1 = $wb.Worksheets.Item("12t")
Remove-StudentRows $ws1 @("2102030089")

$ws2 = $wb.Worksheets.Item("12i")
Remove-StudentRows $ws2 @("2102030089", "2102030073")
